# Apply the latest symbol-list refresh (price/volume updates, plus the
# UpBots/LEO/.../NitroEx row reshuffle) produced by the GitHub Actions
# scraper run on Fri Feb 10 09:31:36 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    # Force the cell to stay plain text so Excel does not silently
    # reinterpret price/percentage-looking strings as numbers.
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "309.24"
Set-TextValue "E2" "-3.91%"
Set-TextValue "D3" "40.24"
Set-TextValue "E3" "-5.94%"
Set-TextValue "D4" "5.127"
Set-TextValue "E4" "-0.42%"
Set-TextValue "D5" "0.07757"
Set-TextValue "E5" "-5.21%"
Set-TextValue "D6" "4.254"
Set-TextValue "E6" "-0.71%"
Set-TextValue "D7" "1.628"
Set-TextValue "E7" "-10.01%"
Set-TextValue "D8" "0.8817"
Set-TextValue "E8" "-5.45%"
Set-TextValue "E9" "-8.02%"
Set-TextValue "D10" "0.1747"
Set-TextValue "E10" "-6.40%"
Set-TextValue "D11" "0.08990"
Set-TextValue "E11" "-5.26%"
Set-TextValue "D12" "0.04424"
Set-TextValue "E12" "-4.83%"
Set-TextValue "E13" "-0.27%"
Set-TextValue "D14" "0.001261"
Set-TextValue "E14" "-2.87%"
Set-TextValue "D15" "0.005829"
Set-TextValue "E15" "1.10%"
Set-TextValue "B16" "UpBots"
Set-TextValue "C16" "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-TextValue "D16" "0.007491"
Set-TextValue "E16" "2,413.46%"
Set-TextValue "B17" "LEO"
Set-TextValue "C17" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D17" "3.354"
Set-TextValue "E17" "-0.32%"
Set-TextValue "B18" "BTSEToken"
Set-TextValue "C18" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D18" "2.420"
Set-TextValue "E18" "-4.28%"
Set-TextValue "B19" "BitpandaEcosystemToken"
Set-TextValue "C19" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D19" "0.3276"
Set-TextValue "E19" "-2.97%"
Set-TextValue "B20" "MCDex"
Set-TextValue "C20" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D20" "7.027"
Set-TextValue "E20" "-5.10%"
Set-TextValue "B21" "ProBitToken"
Set-TextValue "C21" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D21" "0.1339"
Set-TextValue "E21" "-3.58%"
Set-TextValue "B22" "ZBToken"
Set-TextValue "C22" "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue "D22" "0.2787"
Set-TextValue "E22" "10.50%"
Set-TextValue "B23" "CoinExToken"
Set-TextValue "C23" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D23" "0.04184"
Set-TextValue "E23" "0.38%"
Set-TextValue "B24" "BitKan"
Set-TextValue "C24" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D24" "0.001202"
Set-TextValue "E24" "-3.58%"
Set-TextValue "B25" "HotbitToken"
Set-TextValue "C25" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D25" "0.004080"
Set-TextValue "E25" "-6.46%"
Set-TextValue "B26" "NitroEx"
Set-TextValue "C26" "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D26" "0.0001301"
Set-TextValue "E26" "8.34%"
Set-TextValue "D38" "0.02376"
Set-TextValue "E38" "-13.76%"
Set-TextValue "D39" "0.05224"
Set-TextValue "E39" "-6.38%"
Set-TextValue "D40" "0.007940"
Set-TextValue "E40" "-1.52%"
Set-TextValue "D41" "0.1328"
Set-TextValue "E41" "-5.00%"
Set-TextValue "E42" "-3.25%"
Set-TextValue "D43" "0.001958"
Set-TextValue "E43" "-6.48%"
Set-TextValue "D44" "0.008762"
Set-TextValue "E44" "15.89%"
Set-TextValue "D45" "0.3354"
Set-TextValue "E45" "-4.09%"
Set-TextValue "D46" "0.00006549"
Set-TextValue "E46" "-6.09%"
Set-TextValue "E47" "-0.06%"
Set-TextValue "E48" "98.37%"
Set-TextValue "D49" "0.002723"
Set-TextValue "E49" "-21.71%"
Set-TextValue "D50" "0.00002100"
Set-TextValue "E50" "-0.06%"
Set-TextValue "D51" "0.0002000"
Set-TextValue "E51" "-0.06%"
